# "Correccion Alex Preciarios en base a volumetrias" -- fills in the id
# (column A) codes for the category/sub-category header rows that were
# previously left blank, so the "cantidad utilizada" rollups key off a
# real id starting at/accumulating from 0, per the commit message.
#
# The new ids are written in the same order the author introduced them in
# the shared-string table (B02, B03, C, C01, then the A0203.. family) so
# the underlying workbook's shared-string list grows in that same order.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A94").Value = "B02"
$ws.Range("A95").Value = "B03"
$ws.Range("A97").Value = "C"
$ws.Range("A98").Value = "C01"
$ws.Range("A67").Value = "A0203"
$ws.Range("A72").Value = "A03"
$ws.Range("A73").Value = "A0301"
$ws.Range("A76").Value = "A04"
$ws.Range("A77").Value = "A0401"
$ws.Range("A79").Value = "A0402"
$ws.Range("A83").Value = "A0403"

# Restore the view state captured in the saved workbook: scrolled down so
# row 113 is at the top, with B84 selected/active.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 113
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B84").Select()
